# Generate Report for Handback
#
# This mirrors the localization "handback" step: each language sheet
# (zh-cn / de-de) gets its "Latest Target File" + "Latest Handback File"
# columns populated with the handed-back files, the "Latest Handback
# DateTime" stamped with the real completion time, and the Overview
# sheet's per-language Status columns flipped from "Ready for handoff"
# to "Handed back: in sync with en-US".

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)   # Overview
$ws2 = $wb.Worksheets.Item(2)   # zh-cn
$ws3 = $wb.Worksheets.Item(3)   # de-de

$repoBase = "https://github.com/OpenLocalizationTestOrg/oltest/blob/5018b902bf00d845d8c568f51098b7a7b6f681c0/e2e/"

$file1Md   = "15647dbd-6c0d-4ef2-91fe-997abde31ff5.md"
$file2Md   = "44b555cc-9929-4ecf-8cb1-27f61c911242.md"

$file1ZhXlf = "15647dbd-6c0d-4ef2-91fe-997abde31ff5.a58327ac70b79301fbf96170dd6522f691524e82.zh-cn.xlf"
$file2ZhXlf = "44b555cc-9929-4ecf-8cb1-27f61c911242.9c005b14852cb5393981f3a9db9a08e2aeca3331.zh-cn.xlf"
$file1DeXlf = "15647dbd-6c0d-4ef2-91fe-997abde31ff5.a58327ac70b79301fbf96170dd6522f691524e82.de-de.xlf"
$file2DeXlf = "44b555cc-9929-4ecf-8cb1-27f61c911242.9c005b14852cb5393981f3a9db9a08e2aeca3331.de-de.xlf"

$zhHandbackTime = "2016-08-13 02:34:05"
$deHandbackTime = "2016-08-13 02:34:15"

$statusText = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# Overview sheet: flip the per-language status cells.
# ---------------------------------------------------------------------
$ws1.Range("E2").Value = $statusText
$ws1.Range("F2").Value = $statusText
$ws1.Range("E3").Value = $statusText
$ws1.Range("F3").Value = $statusText

$ws1.Columns.Item(5).ColumnWidth = 29.17
$ws1.Columns.Item(6).ColumnWidth = 29.17

# ---------------------------------------------------------------------
# zh-cn sheet: fill in Latest Target File / Latest Handback File /
# Latest Handback DateTime for both rows, and hyperlink the target file
# column the same way column A (Source File Name) is hyperlinked.
# ---------------------------------------------------------------------
$ws2.Hyperlinks.Add($ws2.Range("I2"), ($repoBase + $file1Md), [Type]::Missing, [Type]::Missing, $file1Md)
$ws2.Range("I2").Style = "HyperLink"
$ws2.Range("J2").Value = $file1ZhXlf
$ws2.Range("K2").Value = $zhHandbackTime

$ws2.Hyperlinks.Add($ws2.Range("I3"), ($repoBase + $file2Md), [Type]::Missing, [Type]::Missing, $file2Md)
$ws2.Range("I3").Style = "HyperLink"
$ws2.Range("J3").Value = $file2ZhXlf
$ws2.Range("K3").Value = $zhHandbackTime

$ws2.Columns.Item(3).ColumnWidth = 29.17
$ws2.Columns.Item(9).ColumnWidth = 39.17
$ws2.Columns.Item(10).ColumnWidth = 39.17

# ---------------------------------------------------------------------
# de-de sheet: same shape as zh-cn, different target language files.
# ---------------------------------------------------------------------
$ws3.Hyperlinks.Add($ws3.Range("I2"), ($repoBase + $file1Md), [Type]::Missing, [Type]::Missing, $file1Md)
$ws3.Range("I2").Style = "HyperLink"
$ws3.Range("J2").Value = $file1DeXlf
$ws3.Range("K2").Value = $deHandbackTime

$ws3.Hyperlinks.Add($ws3.Range("I3"), ($repoBase + $file2Md), [Type]::Missing, [Type]::Missing, $file2Md)
$ws3.Range("I3").Style = "HyperLink"
$ws3.Range("J3").Value = $file2DeXlf
$ws3.Range("K3").Value = $deHandbackTime

$ws3.Columns.Item(3).ColumnWidth = 29.17
$ws3.Columns.Item(9).ColumnWidth = 39.17
$ws3.Columns.Item(10).ColumnWidth = 39.17
